$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---

# Cells whose new value could be misread as a plain number: force text
# formatting, assign, then restore the default "Normal" style so no
# stray formatting diff is left behind.
$textForceCells = @("D5", "D6", "D7", "D8", "D13", "D16", "D19", "D21", "D24", "D28", "D29", "D33", "D35", "D36", "D41", "D42", "D44", "D45", "D46", "D49", "D50")
foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D5").Value = "532.58"
$ws.Range("D6").Value = "133.22"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "0.566"
$ws.Range("D13").Value = "0.328"
$ws.Range("D16").Value = "22.20"
$ws.Range("D19").Value = "10.56"
$ws.Range("D21").Value = "319.68"
$ws.Range("D24").Value = "65.91"
$ws.Range("D28").Value = "7.44"
$ws.Range("D29").Value = "172.53"
$ws.Range("D33").Value = "6.26"
$ws.Range("D35").Value = "0.997"
$ws.Range("D36").Value = "18.09"
$ws.Range("D41").Value = "0.804"
$ws.Range("D42").Value = "5.12"
$ws.Range("D44").Value = "274.44"
$ws.Range("D45").Value = "131.25"
$ws.Range("D46").Value = "0.589"
$ws.Range("D49").Value = "0.0217"
$ws.Range("D50").Value = "16.71"

foreach ($cell in $textForceCells) {
    $ws.Range($cell).Style = "Normal"
}

# Cells whose new value is unambiguously non-numeric text already
$ws.Range("D2").Value = "58.672.55"
$ws.Range("D3").Value = "2.488.81"
$ws.Range("D9").Value = "2.499.05"
$ws.Range("D14").Value = "2.922.56"
$ws.Range("D15").Value = "58.632.05"
$ws.Range("D18").Value = "2.502.29"
$ws.Range("D51").Value = "1.747.30"

# --- Volume(1h) (column E) updates ---
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +4.05%  "
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +3.90%  "
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("E45").Value = "  +9.54%  "
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  +2.44%  "

# --- Row 50 / 51 coin identity changes ---
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

